# Align sheet data with "master create db": recode the region column (D)
# from "בני דוד" to the correct region names, and drop the now-obsolete
# duplicate row (row 5), then tidy up the selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "eshcol" (region) column for the three remaining entries.
# Set D3/D4 ("מרכז") before D2 ("דרום") so the new shared-string table
# lands in the same order as the target workbook.
$ws.Range("D3").Value = "מרכז"
$ws.Range("D4").Value = "מרכז"
$ws.Range("D2").Value = "דרום"

# Row 5 duplicated row 4's data (same phone numbers) and is removed.
$ws.Rows(5).Delete()

# Select the (now last) row where the deleted row used to be, matching
# the post-edit selection state saved with the workbook.
$ws.Rows(5).Select() | Out-Null
